$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 24, shifting rows 24-39 down to 25-40.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record's data.
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = "2021-08-09"
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = 100114007
$ws.Cells.Item(24, 7).Value = "Jengibre"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 230
$ws.Cells.Item(24, 11).Value = 13000
$ws.Cells.Item(24, 12).Value = 14000
$ws.Cells.Item(24, 13).Value = 13565
$ws.Cells.Item(24, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(24, 15).Value = "Perú"
$ws.Cells.Item(24, 16).Value = 1043
$ws.Cells.Item(24, 17).Value = 13
$ws.Cells.Item(24, 18).Value = "Hortaliza"
